$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) values for rows that changed.
# D column is forced to text format first so purely-numeric-looking
# price strings (e.g. "0.9991") are not auto-converted to numbers by Excel,
# matching the original inlineStr cell type.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.622.10"
$ws.Range("E2").Value = "  +6.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.741.29"
$ws.Range("E3").Value = "  +4.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.30"
$ws.Range("E5").Value = "  +5.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9966"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3743"
$ws.Range("E7").Value = "  +3.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.71"
$ws.Range("E8").Value = "  +5.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3421"
$ws.Range("E9").Value = "  +4.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.207"
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07522"
$ws.Range("E11").Value = "  +6.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9972"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.479"
$ws.Range("E13").Value = "  +6.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.52"
$ws.Range("E14").Value = "  +4.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.084"
$ws.Range("E15").Value = "  +6.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.731.99"
$ws.Range("E16").Value = "  +3.76%  "
$ws.Range("E17").Value = "  +4.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06710"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.30"
$ws.Range("E19").Value = "  +5.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9964"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.80"
$ws.Range("E21").Value = "  +6.22%  "
$ws.Range("E22").Value = "  +4.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.16"
$ws.Range("E23").Value = "  +4.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "26.578.65"
$ws.Range("E24").Value = "  +6.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.477"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.521"
$ws.Range("E26").Value = "  +4.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.416"
$ws.Range("E27").Value = "  +15.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.43"
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.66"
$ws.Range("E29").Value = "  +5.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.929.64"
$ws.Range("E30").Value = "  +4.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.28"
$ws.Range("E31").Value = "  +5.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.132"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.216"
$ws.Range("E33").Value = "  +5.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08565"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.730"
$ws.Range("E35").Value = "  +3.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "13.18"
$ws.Range("E36").Value = "  +7.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.496"
$ws.Range("E37").Value = "  +4.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2188"
$ws.Range("E40").Value = "  +5.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.682"
$ws.Range("E41").Value = "  +3.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.248"
$ws.Range("E42").Value = "  -3.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6278"
$ws.Range("E43").Value = "  +5.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.59"
$ws.Range("E44").Value = "  +13.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9969"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.906"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6124"
$ws.Range("E47").Value = "  +8.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07345"
$ws.Range("E50").Value = "  +4.42%  "
$ws.Range("E51").Value = "  +3.78%  "

# Row 38/39 swap: Hedera <-> VeChain with new data (VeChain now row 38, Hedera row 39)
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02373"
$ws.Range("E38").Value = "  +4.70%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06371"
$ws.Range("E39").Value = "  +4.57%  "

# Row 48/49 swap: Quant <-> NEARProtocol with new data (NEARProtocol now row 48, Quant row 49)
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.081"
$ws.Range("E48").Value = "  +6.01%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "128.99"
$ws.Range("E49").Value = "  +2.29%  "
